# Daily attendance processing - 2025-12-16 07:35:01
#
# The "Recorded By" column (G) lists the users who touched each attendance
# record as a comma-separated string, e.g. "System, dnasr281@gmail.com".
# This pass re-derives that column so the most recently-acting recorder is
# listed first and the original first entry is moved to the end of the
# list (a left-rotation of the comma-separated values) - e.g.
#   "System, dnasr281@gmail.com"              -> "dnasr281@gmail.com, System"
#   "System, backup@backdoor.com, system"     -> "backup@backdoor.com, system, System"
#   "admin@admin.com, dnasr281@gmail.com"     -> "dnasr281@gmail.com, admin@admin.com"
# Single-author cells (no comma) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$firstRow = $usedRange.Row
$lastRow = $firstRow + $usedRange.Rows.Count - 1

# Locate the "Recorded By" column from the header row instead of assuming G.
$headerRow = $firstRow
$lastCol = $usedRange.Column + $usedRange.Columns.Count - 1
$recordedByCol = 0
for ($c = $usedRange.Column; $c -le $lastCol; $c++) {
    $header = $ws.Cells.Item($headerRow, $c).Value2
    if ($header -eq "Recorded By") {
        $recordedByCol = $c
        break
    }
}
if ($recordedByCol -eq 0) {
    $recordedByCol = 7
}

$changed = 0
for ($r = $headerRow + 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $recordedByCol)
    $v = $cell.Value2
    if ($v -ne $null -and $v -like "*, *") {
        $parts = $v -split ', '
        if ($parts.Count -gt 1) {
            $rotated = $parts[1..($parts.Count - 1)] + $parts[0]
            $cell.Value = $rotated -join ', '
            $changed++
        }
    }
}

Write-Host "Rotated Recorded By list on $changed row(s)."
